$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("F1").Value = "Årsag"
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("H1").Value = "TCV_range"

# Copy header style from existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2..8: new "Årsag" (reason) values in column F, and move the old
# TCV_range value ("120000-140000") into column H
$reasons = @(
    "Anden årsag (angiv hvilken i bemærkninger)",
    "Systemet (uddyb i bemærkninger)",
    "Ikke oplyst",
    "Utilfredshed (Service - uddyb i bemærkninger)",
    "Bruger ikke produktet",
    "Ikke oplyst",
    "Pris"
)

for ($i = 0; $i -lt $reasons.Length; $i++) {
    $row = $i + 2
    $ws.Range("H$row").Value = "120000-140000"
    $ws.Range("F$row").Value = $reasons[$i]
}
